$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D44:D50").Copy()
$ws.Range("Z30").PasteSpecial(-4122)
$ws.Range("Z30:Z36").Merge()
